$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The board's pinout changed (new MCU / new numbering scheme): the "D#"
# Arduino-style pin names are replaced by raw "GPIO#" names, the jumper
# ("JMP1"/"JMP2") pins become real address-select pins, and two more rows
# of pin assignments are appended at the bottom of the table.
# ---------------------------------------------------------------------------

# Column A: IC pin numbers (GPIOxx)
$colA = @{
    3  = "GPIO0"
    4  = "GPIO1"
    5  = "GPIO2"
    6  = "GPIO3"
    7  = "GPIO5"
    8  = "GPIO6"
    9  = "GPIO7"
    10 = "GPIO10"
    11 = "GPIO11"
    12 = "GPIO12"
    13 = "GPIO13"
    14 = "GPIO14"
    15 = "GPIO18"
    16 = "GPIO19"
    17 = "GPIO20"
    18 = "GPIO24"
    19 = "GPIO27"
    20 = "GPIO28"
}

# Column B: peripheral
$colB = @{
    3  = "I2C"
    4  = "I2C"
    5  = ""
    6  = "PIO1"
    7  = ""
    8  = "PIO0"
    9  = "PIO0"
    10 = "PWM5A"
    11 = ""
    12 = ""
    13 = ""
    14 = ""
    15 = "PIO0"
    16 = "PIO0"
    17 = ""
    18 = ""
    19 = ""
    20 = "PWM6A"
}

# Column C: signal
$colC = @{
    3  = "SDA"
    4  = "SCL"
    5  = "M2_DIR"
    6  = "NEOPIXEL"
    7  = "M2_FLAG"
    8  = "ENC2_B"
    9  = "ENC2_A"
    10 = "M2_PWM"
    11 = "M2_ENABLE"
    12 = "ADDR_0"
    13 = "ADDR_1"
    14 = "ADDR_2"
    15 = "ENC1_B"
    16 = "ENC1_A"
    17 = "M1_DIR"
    18 = "M1_FLAG"
    19 = "M1_ENABLE"
    20 = "M1_PWM"
}

# Column D: comment
$colD = @{
    3  = ""
    4  = ""
    5  = ""
    6  = ""
    7  = "Error status flag for driver2"
    8  = ""
    9  = ""
    10 = ""
    11 = ""
    12 = ""
    13 = ""
    14 = ""
    15 = ""
    16 = ""
    17 = ""
    18 = "Error status flag for driver1"
    19 = ""
    20 = ""
}

# Two brand-new rows (19 & 20) are appended at the bottom of the table.
# Seed them from the closest-matching existing bordered rows so they pick
# up the table's border/fill formatting, then overwrite the text.
$ws.Range("A17:D17").Copy($ws.Range("A19:D19"))
$ws.Range("A18:D18").Copy($ws.Range("A20:D20"))

for ($r = 3; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $colA[$r]
    $ws.Cells.Item($r, 2).Value = $colB[$r]
    $ws.Cells.Item($r, 3).Value = $colC[$r]
    $ws.Cells.Item($r, 4).Value = $colD[$r]
}

# Row 15 used to be highlighted (NEOPIXEL/"Not broken out" row) - that
# special shading no longer applies to any row in the new layout, so pull
# its formatting back to the plain bordered look used by the rest of the
# table (copy the format from an already-plain row rather than hand-rolling
# a fresh fill, so we reuse the existing "no fill" style instead of minting
# a new one).
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)

# The jumper-select cells (old C7/C8) used a distinct (unfilled) font
# styling; the new ADDR_* / M2_FLAG / M1_FLAG cells use the regular font,
# so nothing further needs to be done there - new text was already written
# with Cells.Item(...).Value above using the sheet's normal formatting.

# Selection moved off the data table in the saved file.
$ws.Range("G11").Select()
